$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the small table of doubled "hi" strings plus the "Wow!" label,
# matching the shared-strings order from the target workbook:
#   0 hi, 1 hihi, 2 hihihi, 3 hihihihi, 4 Wow!
$ws.Range("A1").Value = "hi"
$ws.Range("A2").Value = "hihi"
$ws.Range("B1").Value = "hihihi"
$ws.Range("B2").Value = "hihihihi"
$ws.Range("C1").Value = "Wow!"

# Leave the selection where the author left it when they saved the file.
$ws.Range("E12").Select() | Out-Null
